$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update D3 status from "Writing" to "Automated" and mark B3 = 1 (Automated test case)
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = "Automated"

# Update the selected cell to B3
$ws.Range("B3").Select()

$wb.Save()
